# Fixed errors with MSM implementation.
#
# This script re-applies the metrics that the analysis tool recomputed:
#  - classFields: the order in which each class's fields are listed was
#    regenerated (reflection order changed), so several rows on the
#    "classFields" sheet need their Field Name / Field Modifier / Field Type
#    swapped to the new order.
#  - classNumberOfLines / methodNumberOfLines: number-of-lines counts for a
#    few constructors / interface methods were fixed from "1" to "0".
#
# Helper: write a TEXT value into a cell (never let Excel reinterpret a
# purely-numeric-looking string such as "0" or "37" as a real number, and
# never let an empty string clear/remove the cell -- the source file keeps
# real, empty shared-string cells). We do this by writing the value with a
# leading apostrophe (forces text, even when empty) and then stripping the
# resulting "quote prefix" cell style by pasting-in the (unformatted) style
# of a pristine helper cell, so the cell ends up as a plain text cell.

function Set-TextValue {
    param(
        $Sheet,
        [string]$CellRef,
        [string]$Text
    )

    $helper = $Sheet.Range("ZZ1000")
    $target = $Sheet.Range($CellRef)

    $target.Formula = "'" + $Text
    $helper.Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null  # xlPasteFormats: restores plain/default style
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# classFields: reorder each class's field rows to match the regenerated
# reflection order.
# ---------------------------------------------------------------------
$classFields = $wb.Worksheets.Item("classFields")

Set-TextValue $classFields "B6"  "factory"
Set-TextValue $classFields "C6"  "private"
Set-TextValue $classFields "D6"  "org.springframework.kafka.core.ConsumerFactory"

Set-TextValue $classFields "B7"  "template"
Set-TextValue $classFields "C7"  "private"
Set-TextValue $classFields "D7"  "org.springframework.kafka.core.KafkaTemplate"

Set-TextValue $classFields "B8"  "LOG"
Set-TextValue $classFields "C8"  "private"
Set-TextValue $classFields "D8"  "org.slf4j.Logger"

Set-TextValue $classFields "B9"  "product"
Set-TextValue $classFields "C9"  ""
Set-TextValue $classFields "D9"  "pl.piomin.stock.domain.Product"

Set-TextValue $classFields "B10" "kafka"
Set-TextValue $classFields "C10" "private"
Set-TextValue $classFields "D10" "org.springframework.kafka.test.EmbeddedKafkaBroker"

Set-TextValue $classFields "B11" "repository"
Set-TextValue $classFields "C11" ""
Set-TextValue $classFields "D11" "pl.piomin.stock.repository.ProductRepository"

Set-TextValue $classFields "B12" "template"
Set-TextValue $classFields "C12" "private"
Set-TextValue $classFields "D12" "org.springframework.kafka.core.KafkaTemplate"

Set-TextValue $classFields "B13" "LOG"
Set-TextValue $classFields "C13" "private"
Set-TextValue $classFields "D13" "org.slf4j.Logger"

Set-TextValue $classFields "B14" "repository"
Set-TextValue $classFields "C14" "private"
Set-TextValue $classFields "D14" "pl.piomin.stock.repository.ProductRepository"

Set-TextValue $classFields "B15" "SOURCE"
Set-TextValue $classFields "C15" "private"
Set-TextValue $classFields "D15" "java.lang.String"

Set-TextValue $classFields "B16" "orderManageService"
Set-TextValue $classFields "C16" ""
Set-TextValue $classFields "D16" "pl.piomin.stock.service.OrderManageService"

Set-TextValue $classFields "B17" "repository"
Set-TextValue $classFields "C17" "private"
Set-TextValue $classFields "D17" "pl.piomin.stock.repository.ProductRepository"

Set-TextValue $classFields "B18" "LOG"
Set-TextValue $classFields "C18" "private"
Set-TextValue $classFields "D18" "org.slf4j.Logger"

# ---------------------------------------------------------------------
# classNumberOfLines: fix number-of-lines metric.
# ---------------------------------------------------------------------
$classLines = $wb.Worksheets.Item("classNumberOfLines")

Set-TextValue $classLines "B4" "0"
Set-TextValue $classLines "B5" "37"
Set-TextValue $classLines "B6" "38"
Set-TextValue $classLines "B8" "23"

# ---------------------------------------------------------------------
# methodNumberOfLines: fix number-of-lines metric (constructors / interface
# methods now report 0 lines instead of 1, other rows simply renumbered).
# ---------------------------------------------------------------------
$methodLines = $wb.Worksheets.Item("methodNumberOfLines")

Set-TextValue $methodLines "C3"  "6"
Set-TextValue $methodLines "C13" "0"
Set-TextValue $methodLines "C15" "0"
Set-TextValue $methodLines "C16" "11"
Set-TextValue $methodLines "C17" "10"
Set-TextValue $methodLines "C18" "10"
Set-TextValue $methodLines "C19" "4"
Set-TextValue $methodLines "C20" "17"
Set-TextValue $methodLines "C21" "13"
Set-TextValue $methodLines "C22" "0"
Set-TextValue $methodLines "C24" "0"
Set-TextValue $methodLines "C26" "9"
Set-TextValue $methodLines "C27" "8"

# ---------------------------------------------------------------------
# Clean up the helper cells used to strip the "quote prefix" style.
# ---------------------------------------------------------------------
$classFields.Range("ZZ1000").Clear() | Out-Null
$classLines.Range("ZZ1000").Clear() | Out-Null
$methodLines.Range("ZZ1000").Clear() | Out-Null
